# Commit: "process refactoring the output"
#
# Renames the Pipette / Buffer / Cassette worksheets to Pipette_1 / Buffer_1 /
# Cassette_1, moves the active tab + selection to the (now 4th) Cassette_1
# sheet, and leaves the Uncut_Sheet_1 tab no longer selected.

$wb = $excel.ActiveWorkbook

# --- Rename the three worksheets -----------------------------------------
$wb.Worksheets.Item("Pipette").Name  = "Pipette_1"
$wb.Worksheets.Item("Buffer").Name   = "Buffer_1"
$wb.Worksheets.Item("Cassette").Name = "Cassette_1"

# --- Re-assert the print areas so the Print_Area defined names track the
#     new worksheet names (Cassette_1!$A$1:$H$42, Pipette_1!$A$1:$G$42,
#     Buffer_1!$A$1:$G$42) --------------------------------------------------
$wb.Worksheets.Item("Pipette_1").PageSetup.PrintArea  = '$A$1:$G$42'
$wb.Worksheets.Item("Buffer_1").PageSetup.PrintArea   = '$A$1:$G$42'
$wb.Worksheets.Item("Cassette_1").PageSetup.PrintArea = '$A$1:$H$42'

# --- Make Cassette_1 (now sheet index 3) the active tab, with C6:H6
#     selected (previously F14 was selected on that sheet; Uncut_Sheet_1
#     loses the tabSelected flag it had before). ---------------------------
$ws = $wb.Worksheets.Item("Cassette_1")
$ws.Activate()
$ws.Range("C6:H6").Select()
